# Update pl_mw (line active power, MW) results for the 380 kV case
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("B2", 0.4526314923384973),
    @("C2", 0.09084829728251975),
    @("D2", 0.2525842607957429),
    @("F2", 1.706025749476368),
    @("G2", 0.9920548021673739),
    @("H2", 1.012311489806052),
    @("I2", 0.7888807512141192),
    @("J2", 0.3509068663016137),
    @("K2", 0.5227810769768269),
    @("N2", 1.710131339581244),
    @("B3", 0.4137551630406051),
    @("C3", 0.08244525588877138),
    @("D3", 0.2438841384559822),
    @("F3", 1.701859145009067),
    @("G3", 0.9910174178869084),
    @("H3", 1.016545079673918),
    @("I3", 0.7931246114590422),
    @("J3", 0.3397061098324201),
    @("K3", 0.4770436605348038),
    @("N3", 1.728768688687948),
    @("B4", 0.3900267168473874),
    @("C4", 0.07731577277471047),
    @("D4", 0.2386586278824581),
    @("F4", 1.700232263775916),
    @("G4", 0.9909909058888502),
    @("H4", 1.01959135829604),
    @("I4", 0.7961333070651975),
    @("J4", 0.3330371592171559),
    @("K4", 0.4491267612859531),
    @("N4", 1.74079249407677),
    @("B5", 0.3803932037354798),
    @("C5", 0.07523302389266462),
    @("D5", 0.2365585387789082),
    @("F5", 1.699803434864336),
    @("G5", 0.9911334231698845),
    @("H5", 1.020945124722274),
    @("I5", 0.7974606442230794),
    @("J5", 0.3303718271877329),
    @("K5", 0.4377925004981762),
    @("N5", 1.74583816882884),
    @("B6", 0.3787957527755452),
    @("C6", 0.0748876417401334),
    @("D6", 0.236211595529582),
    @("F6", 1.69974636733491),
    @("G6", 0.991166343023977),
    @("H6", 1.021176704954968),
    @("I6", 0.7976871623607842),
    @("J6", 0.3299324086623443),
    @("K6", 0.4359130082391687),
    @("N6", 1.746684807874434),
    @("B7", 0.389896649580237),
    @("C7", 0.07728765348547029),
    @("D7", 0.2386301864082014),
    @("F7", 1.700225532527455),
    @("G7", 0.9909922073766211),
    @("H7", 1.019609160600979),
    @("I7", 0.796150798038127),
    @("J7", 0.3330010018738818),
    @("K7", 0.4489737325515648),
    @("N7", 1.740859951257896),
    @("B8", 0.4391977174365991),
    @("C8", 0.0879447008080092),
    @("D8", 0.2495603440777927),
    @("F8", 1.704395732730191),
    @("G8", 0.9915703019076574),
    @("H8", 1.013678502462213),
    @("I8", 0.7902603680555806),
    @("J8", 0.3470015515524096),
    @("K8", 0.5069765422340708),
    @("N8", 1.716436996028463),
    @("B9", 0.5369907806844481),
    @("C9", 0.1090820720106365),
    @("D9", 0.2719158091492488),
    @("F9", 1.719970313126069),
    @("G9", 0.9975572448599337),
    @("H9", 1.005593605247668),
    @("I9", 0.7819094259411017),
    @("J9", 0.3761150330427512),
    @("K9", 0.6220283942263336),
    @("N9", 1.673151201582789),
    @("B10", 0.6095110244877731),
    @("C10", 0.124760056990624),
    @("D10", 0.2889012784538068),
    @("F10", 1.735935168443149),
    @("G10", 1.00493018747693),
    @("H10", 1.001815293049859),
    @("I10", 0.777729667884401),
    @("J10", 0.398525646079861),
    @("K10", 0.7073514779565073),
    @("N10", 1.644160155820639),
    @("B11", 0.6426471350439442),
    @("C11", 0.1319252756347282),
    @("D11", 0.2967501561597032),
    @("F11", 1.744183025649988),
    @("G11", 1.008933703583025),
    @("H11", 1.000566019539903),
    @("I11", 0.7762539363432097),
    @("J11", 0.4089447971371243),
    @("K11", 0.7463395606418146),
    @("N11", 1.631582240852456),
    @("B12", 0.6552156903179878),
    @("C12", 0.1346433537388805),
    @("D12", 0.299739831242789),
    @("F12", 1.747448144139881),
    @("G12", 1.010543378545066),
    @("H12", 1.000160471472554),
    @("I12", 0.7757564070041241),
    @("J12", 0.4129226617303061),
    @("K12", 0.7611281652625337),
    @("N12", 1.626907150874091),
    @("B13", 0.6525079160300606),
    @("C13", 0.1340577552423667),
    @("D13", 0.2990951753124591),
    @("F13", 1.746738632106826),
    @("G13", 1.010192538352698),
    @("H13", 1.000244810134063),
    @("I13", 0.7758608311778659),
    @("J13", 0.4120645169627863),
    @("K13", 0.7579420838221154),
    @("N13", 1.627910103810162),
    @("B14", 0.6436807465349546),
    @("C14", 0.132148798121932),
    @("D14", 0.2969957687683973),
    @("F14", 1.74444880554374),
    @("G14", 1.009064254705777),
    @("H14", 1.000531301492373),
    @("I14", 0.7762117752885871),
    @("J14", 0.4092714094713727),
    @("K14", 0.7475557359379934),
    @("N14", 1.631195856061504),
    @("B15", 0.6382765248334579),
    @("C15", 0.1309801269879927),
    @("D15", 0.2957120950044896),
    @("F15", 1.743064695797997),
    @("G15", 1.008385348632629),
    @("H15", 1.000715579782337),
    @("I15", 0.7764347241618168),
    @("J15", 0.407564766927905),
    @("K15", 0.7411969985524252),
    @("N15", 1.633219924699832),
    @("B16", 0.6073484106397871),
    @("C16", 0.1242924606056022),
    @("D16", 0.2883907857753911),
    @("F16", 1.735415991543292),
    @("G16", 1.004681637999596),
    @("H16", 1.001906384478019),
    @("I16", 0.7778346833615331),
    @("J16", 0.3978492536252816),
    @("K16", 0.7048069872792837),
    @("N16", 1.64499444389935),
    @("B17", 0.5884121925854231),
    @("C17", 0.1201982946547844),
    @("D17", 0.2839306141597149),
    @("F17", 1.730976233691976),
    @("G17", 1.002576053803978),
    @("H17", 1.002757162339023),
    @("I17", 0.7788025978731454),
    @("J17", 0.3919466467679342),
    @("K17", 0.6825272335345289),
    @("N17", 1.652374129511765),
    @("B18", 0.5775343743756025),
    @("C18", 0.1178465666740465),
    @("D18", 0.2813767392240436),
    @("F18", 1.728515341108874),
    @("G18", 1.00142610032303),
    @("H18", 1.003290698037816),
    @("I18", 0.7793993766120408),
    @("J18", 0.3885727448638079),
    @("K18", 0.6697289263527182),
    @("N18", 1.656676161223302),
    @("B19", 0.5738537134992612),
    @("C19", 0.1170508496016396),
    @("D19", 0.2805140183353956),
    @("F19", 1.727698048897707),
    @("G19", 1.001047236780579),
    @("H19", 1.00347893354764),
    @("I19", 0.7796083132836316),
    @("J19", 0.3874340247659944),
    @("K19", 0.6653984709370206),
    @("N19", 1.658142615018644),
    @("B20", 0.5904265606893375),
    @("C20", 0.1206338019660222),
    @("D20", 0.2844042176635071),
    @("F20", 1.73143925437445),
    @("G20", 1.002793869190043),
    @("H20", 1.002662021922944),
    @("I20", 0.7786954150286505),
    @("J20", 0.3925728027596591),
    @("K20", 0.6848972549467192),
    @("N20", 1.65158260440851),
    @("B21", 0.6462729435519634),
    @("C21", 0.1327093758707463),
    @("D21", 0.2976119417413372),
    @("F21", 1.745117532657247),
    @("G21", 1.009393116034573),
    @("H21", 1.000445319339619),
    @("I21", 0.7761070301719357),
    @("J21", 0.4100909339969974),
    @("K21", 0.750605791763121),
    @("N21", 1.630228364447747),
    @("B22", 0.6828918386190139),
    @("C22", 0.1406292317894895),
    @("D22", 0.3063457748612279),
    @("F22", 1.754883861039957),
    @("G22", 1.014251961205588),
    @("H22", 0.9993901677036519),
    @("I22", 0.7747727053041444),
    @("J22", 0.421728692243093),
    @("K22", 0.7936937014990519),
    @("N22", 1.616784511483242),
    @("B23", 0.663336794231725),
    @("C23", 0.1363997186676045),
    @("D23", 0.3016750743474859),
    @("F23", 1.749595690546073),
    @("G23", 1.011608682379702),
    @("H23", 0.9999173043271128),
    @("I23", 0.7754521340335216),
    @("J23", 0.4155001131662175),
    @("K23", 0.7706838618740903),
    @("N23", 1.623912819387359),
    @("B24", 0.589515837223729),
    @("C24", 0.1204369026862935),
    @("D24", 0.2841900693153718),
    @("F24", 1.731229637452302),
    @("G24", 1.002695206180732),
    @("H24", 1.002704896559379),
    @("I24", 0.7787437468268905),
    @("J24", 0.3922896566396616),
    @("K24", 0.683825735158365),
    @("N24", 1.651940268030682),
    @("B25", 0.5104166856908137),
    @("C25", 0.1033379956220983),
    @("D25", 0.2657695021108566),
    @("F25", 1.714963767172222),
    @("G25", 0.9954164263711931),
    @("H25", 1.0074012122315),
    @("I25", 0.7838254867748802),
    @("J25", 0.3680604934990583),
    @("K25", 0.5907640706604411),
    @("N25", 1.684367678305129),
)

foreach ($pair in $updates) {
    $ws.Range($pair[0]).Value2 = $pair[1]
}

